# Update "想去人数" (interest count) figures in F column across sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 106
$ws1.Range("F6").Value = 9416
$ws1.Range("F7").Value = 846
$ws1.Range("F9").Value = 1203
$ws1.Range("F10").Value = 1148
$ws1.Range("F12").Value = 97
$ws1.Range("F14").Value = 261
$ws1.Range("F15").Value = 423
$ws1.Range("F18").Value = 1282

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 13

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 106
$ws4.Range("F4").Value = 13
$ws4.Range("F7").Value = 9416
$ws4.Range("F8").Value = 846
$ws4.Range("F10").Value = 1203
$ws4.Range("F11").Value = 1148
$ws4.Range("F13").Value = 97
$ws4.Range("F15").Value = 261
$ws4.Range("F16").Value = 423
$ws4.Range("F19").Value = 1282
